# Generate Report for Handback
# Updates status/date cells on the Overview, zh-cn and de-de sheets to
# reflect a regenerated handback report (new run timestamps and a
# status change from "ht" to "mt" for the zh-cn/de-de 28575c7e... item).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" column (G) for the first two rows
$overview.Range("G2").Value = "2016-08-24 10:16:36"
$overview.Range("G3").Value = "2016-08-24 10:16:36"

# zh-cn sheet: Priority column (E) changes from "ht" to "mt"
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E3").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H) and Correspond Handback Datetime (K)
$zhcn.Range("H2").Value = "2016-08-24 10:16:31"
$zhcn.Range("H3").Value = "2016-08-24 10:16:31"
$zhcn.Range("K2").Value = "2016-08-24 10:16:48"
$zhcn.Range("K3").Value = "2016-08-24 10:16:48"

# de-de sheet: Priority column (E) changes from "ht" to "mt"
$dede.Range("E2").Value = "mt"
$dede.Range("E3").Value = "mt"

# de-de sheet: Correspond Handoff Datetime (H) mirrors the Overview date
$dede.Range("H2").Value = "2016-08-24 10:16:36"
$dede.Range("H3").Value = "2016-08-24 10:16:36"

# de-de sheet: Correspond Handback Datetime (K)
$dede.Range("K2").Value = "2016-08-24 10:16:55"
$dede.Range("K3").Value = "2016-08-24 10:16:55"
